$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

# Rename table column "Coluna1" -> "23/out" (this updates the header cell B2 and the
# sharedStrings / table column name).
$ws.Range("B2").Value = "23/out"

# The structured-reference formula in the "Total" row needs to follow the new
# column name, same as Excel auto-updates it when a table header is renamed.
$ws.Range("B9").Formula = "=SUBTOTAL(109,Tabela2[23/out])"

# Fill in the price values that were typed into the table body.
$ws.Range("B3").Value = 1054
$ws.Range("B4").Value = 600
$ws.Range("B5").Value = 323
$ws.Range("B6").Value = 297
$ws.Range("B7").Value = 480
$ws.Range("B8").Value = 520

# Selection ends on B3 as the last active cell.
$ws.Range("B3").Select()

$wb.Save()
